# Updates the cryptos list (prices / 1h volume %, plus two row re-orderings)
# per the Fri Sep 13 22:55:36 UTC 2024 GitHub Actions refresh.
#
# Columns D/E hold plain text (prices like "60.579.78", percentages like
# "  +4.20%  "). Some of the new D-column prices (e.g. "557.43") look like
# ordinary decimal numbers, so Excel's COM layer would otherwise silently
# reinterpret them as numeric values. To keep them as text we briefly force
# the cell to Text format ("@"), assign the value, then restore the cell's
# original ("Normal") style so no visible formatting change is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.579.78"
$ws.Range("E2").Value = "  +4.20%  "
$ws.Range("D3").Value = "2.450.96"
$ws.Range("E3").Value = "  +3.80%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.28%  "
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("E8").Value = "  +1.81%  "
$ws.Range("E9").Value = "  +4.23%  "
$ws.Range("E10").Value = "  +4.58%  "
$ws.Range("E11").Value = "  +1.79%  "
$ws.Range("E12").Value = "  -2.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.08"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.57%  "
$ws.Range("D14").Value = "2.883.36"
$ws.Range("E14").Value = "  +3.51%  "
$ws.Range("D15").Value = "60.467.56"
$ws.Range("E15").Value = "  +3.92%  "
$ws.Range("E16").Value = "  +4.41%  "
$ws.Range("D17").Value = "2.451.83"
$ws.Range("E17").Value = "  +3.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +7.16%  "
$ws.Range("E19").Value = "  +3.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "336.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.93%  "
$ws.Range("E21").Value = "  +2.19%  "
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.04%  "
$ws.Range("E24").Value = "  +2.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("E27").Value = "  -0.54%  "
$ws.Range("E28").Value = "  +7.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.81"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.46%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.66%  "
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.32"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.86"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.82%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("E35").Value = "  +5.47%  "
$ws.Range("E36").Value = "  +0.82%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("E38").Value = "  +0.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "40.17"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.418"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "318.42"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.34%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "144.62"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.84%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.74"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.05"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.04%  "
$ws.Range("E45").Value = "  +1.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0526"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.576"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.409"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.26%  "
$ws.Range("E49").Value = "  +2.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.23%  "
$ws.Range("E51").Value = "  +6.16%  "
